$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# CU-13 (Consultar historial de pago de profesores) - row 17
$ws.Range("E17").Value = "planificado"
$ws.Range("F17").Value = 1

# CU-14 (Registrar pago de profesor) - row 18
$ws.Range("E18").Value = "planificado"
$ws.Range("F18").Value = 1
